# Applies odds/stat corrections to Sheet1 for the FlashScore weekly games export.
# Each assignment corresponds to one changed <c> cell value in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 7).Value = 2.38  # G2: 2.3 -> 2.38
$ws.Cells.Item(2, 8).Value = 3  # H2: 3.1 -> 3
$ws.Cells.Item(2, 11).Value = 7.5  # K2: 8 -> 7.5
$ws.Cells.Item(2, 14).Value = 2.3  # N2: 2.25 -> 2.3
$ws.Cells.Item(2, 15).Value = 1.6  # O2: 1.62 -> 1.6
$ws.Cells.Item(2, 18).Value = 1.95  # R2: 1.91 -> 1.95
$ws.Cells.Item(2, 19).Value = 1.8  # S2: 1.91 -> 1.8
$ws.Cells.Item(2, 26).Value = 7.5  # Z2: 8 -> 7.5
$ws.Cells.Item(2, 27).Value = 5.5  # AA2: 6 -> 5.5
$ws.Cells.Item(2, 30).Value = 351  # AD2: 301 -> 351

# Row 3
$ws.Cells.Item(3, 7).Value = 1.76  # G3: 1.7 -> 1.76
$ws.Cells.Item(3, 8).Value = 3.25  # H3: 3.4 -> 3.25
$ws.Cells.Item(3, 9).Value = 5  # I3: 5.5 -> 5
$ws.Cells.Item(3, 10).Value = 1.11  # J3: 1.1 -> 1.11
$ws.Cells.Item(3, 11).Value = 6.5  # K3: 7 -> 6.5
$ws.Cells.Item(3, 12).Value = 1.53  # L3: 1.44 -> 1.53
$ws.Cells.Item(3, 13).Value = 2.38  # M3: 2.63 -> 2.38
$ws.Cells.Item(3, 14).Value = 2.7  # N3: 2.4 -> 2.7
$ws.Cells.Item(3, 15).Value = 1.44  # O3: 1.53 -> 1.44
$ws.Cells.Item(3, 16).Value = 1.57  # P3: 1.53 -> 1.57
$ws.Cells.Item(3, 17).Value = 2.25  # Q3: 2.38 -> 2.25
$ws.Cells.Item(3, 21).Value = 7  # U3: 6.5 -> 7
$ws.Cells.Item(3, 22).Value = 9.5  # V3: 9 -> 9.5
$ws.Cells.Item(3, 23).Value = 13  # W3: 12 -> 13
$ws.Cells.Item(3, 24).Value = 19  # X3: 17 -> 19
$ws.Cells.Item(3, 26).Value = 6.5  # Z3: 7 -> 6.5
$ws.Cells.Item(3, 27).Value = 6.5  # AA3: 7 -> 6.5
$ws.Cells.Item(3, 29).Value = 101  # AC3: 81 -> 101
$ws.Cells.Item(3, 31).Value = 9.5  # AE3: 11 -> 9.5
$ws.Cells.Item(3, 32).Value = 23  # AF3: 26 -> 23
$ws.Cells.Item(3, 34).Value = 51  # AH3: 67 -> 51

# Row 5
$ws.Cells.Item(5, 7).Value = 2  # G5: 2.05 -> 2
$ws.Cells.Item(5, 8).Value = 3.1  # H5: 3.2 -> 3.1
$ws.Cells.Item(5, 9).Value = 4.2  # I5: 3.8 -> 4.2
$ws.Cells.Item(5, 11).Value = 7.5  # K5: 8 -> 7.5
$ws.Cells.Item(5, 12).Value = 1.44  # L5: 1.4 -> 1.44
$ws.Cells.Item(5, 13).Value = 2.63  # M5: 2.75 -> 2.63
$ws.Cells.Item(5, 14).Value = 2.35  # N5: 2.3 -> 2.35
$ws.Cells.Item(5, 15).Value = 1.57  # O5: 1.6 -> 1.57
$ws.Cells.Item(5, 16).Value = 1.53  # P5: 1.5 -> 1.53
$ws.Cells.Item(5, 17).Value = 2.38  # Q5: 2.5 -> 2.38
$ws.Cells.Item(5, 18).Value = 2.1  # R5: 2 -> 2.1
$ws.Cells.Item(5, 19).Value = 1.67  # S5: 1.73 -> 1.67
$ws.Cells.Item(5, 20).Value = 6  # T5: 6.5 -> 6
$ws.Cells.Item(5, 21).Value = 8.5  # U5: 9 -> 8.5
$ws.Cells.Item(5, 22).Value = 9.5  # V5: 9 -> 9.5
$ws.Cells.Item(5, 27).Value = 6  # AA5: 6.5 -> 6
$ws.Cells.Item(5, 28).Value = 19  # AB5: 17 -> 19
$ws.Cells.Item(5, 31).Value = 9  # AE5: 9.5 -> 9
$ws.Cells.Item(5, 33).Value = 15  # AG5: 13 -> 15
$ws.Cells.Item(5, 35).Value = 41  # AI5: 34 -> 41

# Row 6
$ws.Cells.Item(6, 10).Value = 1.05  # J6: 1.06 -> 1.05
$ws.Cells.Item(6, 11).Value = 11  # K6: 10 -> 11
$ws.Cells.Item(6, 12).Value = 1.25  # L6: 1.29 -> 1.25
$ws.Cells.Item(6, 13).Value = 3.75  # M6: 3.5 -> 3.75
$ws.Cells.Item(6, 14).Value = 1.9  # N6: 1.93 -> 1.9
$ws.Cells.Item(6, 15).Value = 1.95  # O6: 1.93 -> 1.95

# Row 7
$ws.Cells.Item(7, 7).Value = 2.6  # G7: 2.63 -> 2.6
$ws.Cells.Item(7, 9).Value = 3  # I7: 2.9 -> 3
$ws.Cells.Item(7, 10).Value = 1.11  # J7: 1.13 -> 1.11
$ws.Cells.Item(7, 11).Value = 6.5  # K7: 6 -> 6.5
$ws.Cells.Item(7, 12).Value = 1.53  # L7: 1.62 -> 1.53
$ws.Cells.Item(7, 13).Value = 2.38  # M7: 2.2 -> 2.38
$ws.Cells.Item(7, 14).Value = 2.7  # N7: 2.88 -> 2.7
$ws.Cells.Item(7, 15).Value = 1.44  # O7: 1.4 -> 1.44
$ws.Cells.Item(7, 18).Value = 2.2  # R7: 2.25 -> 2.2
$ws.Cells.Item(7, 19).Value = 1.62  # S7: 1.57 -> 1.62
$ws.Cells.Item(7, 24).Value = 26  # X7: 29 -> 26
$ws.Cells.Item(7, 26).Value = 6  # Z7: 5.5 -> 6
$ws.Cells.Item(7, 32).Value = 13  # AF7: 12 -> 13
$ws.Cells.Item(7, 34).Value = 34  # AH7: 29 -> 34

# Row 9
$ws.Cells.Item(9, 14).Value = 1.95  # N9: 1.98 -> 1.95
$ws.Cells.Item(9, 15).Value = 1.9  # O9: 1.88 -> 1.9

# Row 10
$ws.Cells.Item(10, 7).Value = 1.38  # G10: 1.42 -> 1.38
$ws.Cells.Item(10, 8).Value = 4.2  # H10: 4 -> 4.2
$ws.Cells.Item(10, 9).Value = 10  # I10: 9.5 -> 10
$ws.Cells.Item(10, 11).Value = 8  # K10: 7.5 -> 8
$ws.Cells.Item(10, 18).Value = 2.75  # R10: 2.63 -> 2.75
$ws.Cells.Item(10, 19).Value = 1.4  # S10: 1.44 -> 1.4
$ws.Cells.Item(10, 21).Value = 5  # U10: 5.5 -> 5
$ws.Cells.Item(10, 23).Value = 8  # W10: 8.5 -> 8
$ws.Cells.Item(10, 27).Value = 8.5  # AA10: 8 -> 8.5
$ws.Cells.Item(10, 29).Value = 126  # AC10: 101 -> 126
$ws.Cells.Item(10, 36).Value = 101  # AJ10: 81 -> 101

# Row 11
$ws.Cells.Item(11, 11).Value = 8  # K11: 7.5 -> 8
$ws.Cells.Item(11, 14).Value = 2.35  # N11: 2.3 -> 2.35
$ws.Cells.Item(11, 15).Value = 1.57  # O11: 1.6 -> 1.57

# Row 18
$ws.Cells.Item(18, 10).Value = 1.06  # J18: 1.07 -> 1.06
$ws.Cells.Item(18, 11).Value = 10  # K18: 9 -> 10
$ws.Cells.Item(18, 12).Value = 1.3  # L18: 1.33 -> 1.3
$ws.Cells.Item(18, 13).Value = 3.4  # M18: 3.25 -> 3.4
$ws.Cells.Item(18, 14).Value = 2.05  # N18: 2.08 -> 2.05
$ws.Cells.Item(18, 15).Value = 1.75  # O18: 1.73 -> 1.75

# Row 23
$ws.Cells.Item(23, 7).Value = 2.12  # G23: 2.18 -> 2.12
$ws.Cells.Item(23, 8).Value = 2.82  # H23: 3.1 -> 2.82
$ws.Cells.Item(23, 9).Value = 3.75  # I23: 3.25 -> 3.75
$ws.Cells.Item(23, 12).Value = 1.47  # L23: 1.45 -> 1.47
$ws.Cells.Item(23, 13).Value = 2.35  # M23: 2.37 -> 2.35
$ws.Cells.Item(23, 14).Value = 2.32  # N23: 2.27 -> 2.32
$ws.Cells.Item(23, 15).Value = 1.47  # O23: 1.5 -> 1.47
$ws.Cells.Item(23, 18).Value = 1.98  # R23: 2 -> 1.98
$ws.Cells.Item(23, 20).Value = 5.8  # T23: 5.9 -> 5.8
$ws.Cells.Item(23, 22).Value = 9  # V23: 9.5 -> 9
$ws.Cells.Item(23, 24).Value = 20  # X23: 21 -> 20
$ws.Cells.Item(23, 25).Value = 37  # Y23: 40 -> 37
$ws.Cells.Item(23, 26).Value = 6.4  # Z23: 7 -> 6.4
$ws.Cells.Item(23, 27).Value = 5.6  # AA23: 6.1 -> 5.6
$ws.Cells.Item(23, 28).Value = 16.5  # AB23: 18 -> 16.5
$ws.Cells.Item(23, 29).Value = 100  # AC23: 110 -> 100
$ws.Cells.Item(23, 31).Value = 8.25  # AE23: 7.8 -> 8.25
$ws.Cells.Item(23, 32).Value = 19  # AF23: 15.5 -> 19
$ws.Cells.Item(23, 33).Value = 13  # AG23: 12 -> 13
$ws.Cells.Item(23, 34).Value = 60  # AH23: 45 -> 60
$ws.Cells.Item(23, 35).Value = 45  # AI23: 35 -> 45
$ws.Cells.Item(23, 36).Value = 55  # AJ23: 50 -> 55
